$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2218934911242604
$ws.Range("C2").Value = 0.5029585798816568
$ws.Range("J2").Value = 0.01775147928994083
$ws.Range("P2").Value = 0.1331360946745562
$ws.Range("S2").Value = 0.1242603550295858
$ws.Range("B3").Value = 0.01734104046242774
$ws.Range("C3").Value = 0.005780346820809248
$ws.Range("J3").Value = 0.01734104046242774
$ws.Range("P3").Value = 0.6936416184971098
$ws.Range("S3").Value = 0.2658959537572254
$ws.Range("J4").Value = 0.01754385964912281
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3157894736842105
$ws.Range("B6").Value = 0.04411764705882353
$ws.Range("D6").Value = 0.009803921568627451
$ws.Range("F6").Value = 0.05392156862745098
$ws.Range("J6").Value = 0.2794117647058824
$ws.Range("O6").Value = 0.02941176470588235
$ws.Range("Q6").Value = 0.1470588235294118
$ws.Range("R6").Value = 0.08333333333333333
$ws.Range("S6").Value = 0.3529411764705883
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.02116402116402116
$ws.Range("F7").Value = 0.05291005291005291
$ws.Range("J7").Value = 0.1164021164021164
$ws.Range("O7").Value = 0.01587301587301587
$ws.Range("Q7").Value = 0.164021164021164
$ws.Range("R7").Value = 0.06878306878306878
$ws.Range("S7").Value = 0.4497354497354497
$ws.Range("B8").Value = 0.07847082494969819
$ws.Range("D8").Value = 0.01810865191146881
$ws.Range("F8").Value = 0.05835010060362173
$ws.Range("J8").Value = 0.08853118712273642
$ws.Range("O8").Value = 0.008048289738430584
$ws.Range("Q8").Value = 0.1468812877263581
$ws.Range("R8").Value = 0.096579476861167
$ws.Range("S8").Value = 0.5050301810865191
$ws.Range("B9").Value = 0.09767441860465116
$ws.Range("F9").Value = 0.06511627906976744
$ws.Range("J9").Value = 0.07906976744186046
$ws.Range("O9").Value = 0.0186046511627907
$ws.Range("Q9").Value = 0.1767441860465116
$ws.Range("R9").Value = 0.07906976744186046
$ws.Range("S9").Value = 0.4837209302325581
$ws.Range("B10").Value = 0.1271763815291446
$ws.Range("D10").Value = 0.03330809992429978
$ws.Range("F10").Value = 0.05526116578349735
$ws.Range("J10").Value = 0.09386828160484481
$ws.Range("O10").Value = 0.01362604087812263
$ws.Range("Q10").Value = 0.2172596517789553
$ws.Range("R10").Value = 0.07494322482967448
$ws.Range("S10").Value = 0.384557153671461
$ws.Range("G11").Value = 0.1331058020477816
$ws.Range("J11").Value = 0.10580204778157
$ws.Range("K11").Value = 0.1774744027303754
$ws.Range("L11").Value = 0.5699658703071673
$ws.Range("S11").Value = 0.0136518771331058
$ws.Range("G12").Value = 0.6949152542372882
$ws.Range("J12").Value = 0.2033898305084746
$ws.Range("K12").Value = 0.01694915254237288
$ws.Range("L12").Value = 0.04519774011299435
$ws.Range("S12").Value = 0.03954802259887006
$ws.Range("F13").Value = 0.02272727272727273
$ws.Range("G13").Value = 0.5909090909090909
$ws.Range("J13").Value = 0.3181818181818182
$ws.Range("S13").Value = 0.06818181818181818
$ws.Range("F15").Value = 0.01293103448275862
$ws.Range("H15").Value = 0.1681034482758621
$ws.Range("I15").Value = 0.05603448275862069
$ws.Range("J15").Value = 0.3706896551724138
$ws.Range("K15").Value = 0.06896551724137931
$ws.Range("O15").Value = 0.08620689655172414
$ws.Range("S15").Value = 0.2370689655172414
$ws.Range("F16").Value = 0.03
$ws.Range("H16").Value = 0.195
$ws.Range("I16").Value = 0.08500000000000001
$ws.Range("J16").Value = 0.385
$ws.Range("K16").Value = 0.13
$ws.Range("M16").Value = 0.005
$ws.Range("O16").Value = 0.065
$ws.Range("S16").Value = 0.105
$ws.Range("F17").Value = 0.0196078431372549
$ws.Range("H17").Value = 0.1612200435729848
$ws.Range("I17").Value = 0.1132897603485839
$ws.Range("J17").Value = 0.4139433551198257
$ws.Range("K17").Value = 0.08278867102396514
$ws.Range("M17").Value = 0.02178649237472767
$ws.Range("O17").Value = 0.06535947712418301
$ws.Range("S17").Value = 0.1220043572984749
$ws.Range("F18").Value = 0.01036269430051814
$ws.Range("H18").Value = 0.1813471502590674
$ws.Range("I18").Value = 0.09326424870466321
$ws.Range("J18").Value = 0.4196891191709844
$ws.Range("K18").Value = 0.09326424870466321
$ws.Range("M18").Value = 0.01036269430051814
$ws.Range("O18").Value = 0.06217616580310881
$ws.Range("S18").Value = 0.1295336787564767
$ws.Range("F19").Value = 0.0110803324099723
$ws.Range("H19").Value = 0.2174515235457064
$ws.Range("I19").Value = 0.0796398891966759
$ws.Range("J19").Value = 0.3857340720221606
$ws.Range("K19").Value = 0.09626038781163435
$ws.Range("M19").Value = 0.0221606648199446
$ws.Range("N19").Value = 0.001385041551246537
$ws.Range("O19").Value = 0.06578947368421052
$ws.Range("S19").Value = 0.1204986149584488
